# product_import_template.xlsx edit
# 1) Rename the first sheet "Dữ liệu mẫu" -> "Dữ liệu sản phẩm"
# 2) Adjust the "Lưu ý quan trọng" row on the "Hướng dẫn" sheet so its
#    height goes from 26.4 to 13.2 (row 7).

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item(1)
$wsGuide = $wb.Worksheets.Item(2)

# Rename the sample-data sheet to reflect that it now holds product data.
$wsData.Name = "Dữ liệu sản phẩm"

# Row 7 ("Lưu ý quan trọng") no longer needs the taller two-line height.
$wsGuide.Rows.Item(7).RowHeight = 13.2
